$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill rows 3 through 11 with the "prueba" test data, mirroring row 2's
# pattern but incrementing the date, document number, materia and
# antecedentes for each row.
$startDate = 44829
$startDoc = 3001

for ($i = 0; $i -lt 9; $i++) {
    $row = 3 + $i
    $n = $i + 2

    $ws.Cells.Item($row, 1).Value = $startDate + $i
    $ws.Cells.Item($row, 2).Value = $startDoc + $i
    $ws.Cells.Item($row, 3).Value = "Carta"
    $ws.Cells.Item($row, 4).Value = "copazo"
    $ws.Cells.Item($row, 5).Value = "lhernandez"
    $ws.Cells.Item($row, 6).Value = "prueba materia $n"
    $ws.Cells.Item($row, 7).Value = "antecedentes prueba $n"
    $ws.Cells.Item($row, 8).Value = "of-0009-21.pdf"
    $ws.Cells.Item($row, 9).Value = "gbenavides"
}

$ws.Range("A2:I11").Select()
$excel.ActiveWindow.RangeSelection.Item(1,1).Activate() | Out-Null
$ws.Range("I2").Activate()
